# "add conference record by lhw"
#
# Adds a second meeting record (row 3) to the discussion-log sheet:
#   A3 = 2015-03-06 18:30  (one day after the existing A2 record)
#   B3 = 刘瀚文 (same host as row 2 -> reuses shared string 4)
#   C3 = long description of the new meeting's content
#   D3 = "80分钟"
# and widens/wraps column C so the new (and existing) descriptions
# read nicely, which reflows the row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row 3 -------------------------------------------------------
# Write D3 before C3 so the shared-string table picks up "80分钟"
# (index 7) ahead of the long description (index 8), matching the
# order new strings were introduced in the sheet.
$ws.Range("A3").Value = 42069.770833333336
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = "刘瀚文"
$ws.Range("D3").Value = "80分钟"
$ws.Range("C3").Value = "讨论作业要求，确定成员分工，制定计划，分配具体文档、编码工作。介绍github、Travis-ci、word的使用方法"

# --- column C: wider + word-wrapped (also applies to header/row 2) --
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Range("C1:C3").WrapText = $true

# --- row heights reflowed by the wider/wrapped column C -------------
$ws.Rows.Item(2).RowHeight = 40.5
$ws.Rows.Item(3).RowHeight = 54

# --- leftover cursor/selection state from the editing session -------
$ws.Range("D6:D9").Select()
